$wb = $excel.ActiveWorkbook

# --- 1. Summary sheet ("总计"): shift current 2022-Q2 row down and add new 2022-Q3 row ---
$wsTotal = $wb.Worksheets.Item("总计")

# Duplicate the formatting of row 2 into row 3 (keeps the "s=2" style on column A)
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Row 3 becomes what row 2 used to hold (the 2022-Q2 totals)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.44

# Row 2 becomes the new 2022-Q3 totals
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.37

# --- 2. Duplicate the existing "2022-Q2" detail sheet so the old data is preserved ---
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsQ2)

# The original sheet (still holding the old Q2 figures) becomes the new Q3 sheet ...
$wsQ2.Name = "2022-Q3"

# ... and the freshly made copy keeps the old data, renamed back to 2022-Q2
$wsCopy = $wb.Worksheets.Item("2022-Q2 (2)")
$wsCopy.Name = "2022-Q2"

# --- 3. Overwrite the (renamed) 2022-Q3 sheet with the new quarter's fund data ---
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

$wsQ3.Range("B2:B3").NumberFormat = "@"
$wsQ3.Range("D2:G3").NumberFormat = "@"

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "014734"
$wsQ3.Range("C2").Value = "广发睿合混合A"
$wsQ3.Range("D2").Value = "5.96"
$wsQ3.Range("E2").Value = "86.96"
$wsQ3.Range("F2").Value = "4.98"
$wsQ3.Range("G2").Value = "0.2968"
$wsQ3.Range("H2").Value = 9

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "014735"
$wsQ3.Range("C3").Value = "广发睿合混合C"
$wsQ3.Range("D3").Value = "1.47"
$wsQ3.Range("E3").Value = "86.96"
$wsQ3.Range("F3").Value = "4.98"
$wsQ3.Range("G3").Value = "0.0732"
$wsQ3.Range("H3").Value = 9
